$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: C11/D11 status changes from "No" (Bad style) to "Done" (Good style)
$ws.Range("C11").Value = "Done"
$ws.Range("C2:D2").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)

# New row 21: DB change / Add skills to the graduate table
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("A21").Value = "DB change"
$ws.Range("B21").Value = "Add skills to the graduate table"

$excel.CutCopyMode = $false
